{"js": "// The document was edited so that three paragraphs, which previously had\n// their text split across many single-word/single-space runs, now each\n// carry their full text inside one run. Rebuild each target paragraph's\n// text (by style) as a single contiguous run, without altering the text\n// content itself.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// style -> full merged text (exactly what the many small runs already\n// concatenate to, so the visible content is unchanged).\nconst targets = {\n  \"Title\": \"Answers: Trigonometric identities (radians)\",\n  \"Author\": \"Dzhemma Ruseva\",\n  \"Abstract\": \"A selection of questions on trigonometric identities, using radians to measure angles.\"\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const style = para.style;\n  if (Object.prototype.hasOwnProperty.call(targets, style)) {\n    // Replace the whole paragraph's content with a single run containing\n    // the same text, collapsing the previous run-per-word/space split.\n    para.insertText(targets[style], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document was edited so that three paragraphs, which previously had\n# their text split across many single-word/single-space runs, now each\n# carry their full text inside one run. Rebuild each target paragraph's\n# text (by style) as a single contiguous run, without altering the text\n# content itself.\n#\n# A same-text Range.Text assignment is a silent no-op in this host, so we\n# instead drive it through Find/Replace (scoped to each paragraph's own\n# Range, so we never touch look-alike text elsewhere in the document),\n# which does rebuild the run(s) that back the paragraph's text.\n\n$d = $word.ActiveDocument\n\n$targets = @{\n    \"Title\"    = \"Answers: Trigonometric identities (radians)\"\n    \"Author\"   = \"Dzhemma Ruseva\"\n    \"Abstract\" = \"A selection of questions on trigonometric identities, using radians to measure angles.\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Range.Style.NameLocal\n    if ($targets.ContainsKey($styleName)) {\n        $text = $targets[$styleName]\n        $rng = $p.Range\n        $rng.Find.ClearFormatting()\n        $rng.Find.Replacement.ClearFormatting()\n        $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null\n    }\n}\n"}
